$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.800.02'
$ws.Range("E2").Value = '  -4.64%  '
$ws.Range("D3").Value = '2.455.96'
$ws.Range("E3").Value = '  -5.95%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''544.74'
$ws.Range("E5").Value = '  -5.17%  '
$ws.Range("D6").Value = '''144.75'
$ws.Range("E6").Value = '  -7.58%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''0.597'
$ws.Range("E8").Value = '  -4.36%  '
$ws.Range("D9").Value = '2.456.04'
$ws.Range("E9").Value = '  -5.85%  '
$ws.Range("D10").Value = '''0.106'
$ws.Range("E10").Value = '  -10.24%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '''5.33'
$ws.Range("E12").Value = '  -8.60%  '
$ws.Range("D13").Value = '''0.350'
$ws.Range("E13").Value = '  -8.10%  '
$ws.Range("D14").Value = '''25.77'
$ws.Range("E14").Value = '  -8.64%  '
$ws.Range("D15").Value = '2.900.14'
$ws.Range("E15").Value = '  -5.92%  '
$ws.Range("D16").Value = '''0.0000162'
$ws.Range("E16").Value = '  -10.30%  '
$ws.Range("D17").Value = '60.731.13'
$ws.Range("E17").Value = '  -4.57%  '
$ws.Range("D18").Value = '2.460.14'
$ws.Range("E18").Value = '  -6.65%  '
$ws.Range("D19").Value = '''11.01'
$ws.Range("E19").Value = '  -8.52%  '
$ws.Range("D20").Value = '''6.90'
$ws.Range("D21").Value = '''4.15'
$ws.Range("E21").Value = '  -8.58%  '
$ws.Range("D22").Value = '''316.26'
$ws.Range("E22").Value = '  -7.99%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '''62.88'
$ws.Range("E24").Value = '  -6.63%  '
$ws.Range("D25").Value = '''1.72'
$ws.Range("E25").Value = '  -5.49%  '
$ws.Range("D26").Value = '0.0₃0968'
$ws.Range("E26").Value = '  -10.94%  '
$ws.Range("D27").Value = '2.576.71'
$ws.Range("E27").Value = '  -5.77%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").Value = '''1.49'
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("D30").Value = '''532.48'
$ws.Range("E30").Value = '  -10.96%  '
$ws.Range("D31").Value = '''8.25'
$ws.Range("E31").Value = '  -10.11%  '
$ws.Range("D32").Value = '''7.62'
$ws.Range("E32").Value = '  -3.68%  '
$ws.Range("E33").Value = '  -8.33%  '
$ws.Range("D34").Value = '''1.88'
$ws.Range("E34").Value = '  -8.95%  '
$ws.Range("E35").Value = '  -10.54%  '
$ws.Range("D36").Value = '''5.84'
$ws.Range("E36").Value = '  -11.74%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '''4.83'
$ws.Range("E38").Value = '  -10.52%  '
$ws.Range("E39").Value = '  -6.92%  '
$ws.Range("D40").Value = '''18.22'
$ws.Range("E40").Value = '  -7.79%  '
$ws.Range("D41").Value = '''144.46'
$ws.Range("E41").Value = '  -6.63%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").Value = '''1.69'
$ws.Range("E43").Value = '  -9.89%  '
$ws.Range("D44").Value = '''39.91'
$ws.Range("E44").Value = '  -3.89%  '
$ws.Range("D45").Value = '''2.27'
$ws.Range("E45").Value = '  -10.61%  '
$ws.Range("D46").Value = '''146.05'
$ws.Range("E46").Value = '  -7.25%  '
$ws.Range("D47").Value = '''3.55'
$ws.Range("E47").Value = '  -9.27%  '
$ws.Range("D48").Value = '''20.71'
$ws.Range("E48").Value = '  -13.14%  '
$ws.Range("D49").Value = '''0.0527'
$ws.Range("E49").Value = '  -10.81%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.580'
$ws.Range("E50").Value = '  -7.86%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.0937'
$ws.Range("E51").Value = '  -6.47%  '
